# "Made changes to grid connection and sun profile"
# Updates the SOLAR annual production figures (column C, "grid connection")
# on the Yearly sheet, and refreshes the sheet's view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yearly")

# Updated SOLAR - Annual prod MWh, MW values for 2020-2024
$ws.Range("C2").Value = 880
$ws.Range("C3").Value = 850
$ws.Range("C4").Value = 700
$ws.Range("C5").Value = 700
$ws.Range("C6").Value = 1000

# Refresh the view: bring it back to the top-left (A1) scroll position,
# reset zoom to 100%, and move the selection to C3
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("C3").Select()
